$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"); copy H1 first so the new
# headers inherit the same bold/centered/bordered style used by the others,
# then overwrite their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 values for rows 2-66
$i0 = @(
    5,
    10,
    7,
    11,
    7,
    7,
    6,
    8,
    7,
    7,
    8,
    8,
    7,
    8,
    8,
    7,
    7,
    7,
    11,
    8,
    7,
    7,
    7,
    8,
    8,
    8,
    7,
    8,
    7,
    6,
    7,
    6,
    7,
    9,
    6,
    9,
    9,
    6,
    6,
    7,
    7,
    7,
    7,
    7,
    6,
    6,
    8,
    7,
    7,
    7,
    7,
    8,
    7,
    7,
    7,
    8,
    8,
    9,
    7,
    7,
    9,
    8,
    8,
    7,
    8
)

# IF values for rows 2-66
$iF = @(
    5,
    10,
    7,
    11,
    7,
    7,
    6,
    8,
    7,
    7,
    8,
    8,
    7,
    8,
    8,
    7,
    7,
    7,
    11,
    8,
    8,
    7,
    7,
    8,
    8,
    8,
    8,
    8,
    7,
    6,
    7,
    7,
    8,
    9,
    6,
    9,
    9,
    6,
    6,
    7,
    7,
    7,
    7,
    7,
    6,
    6,
    8,
    7,
    7,
    7,
    7,
    8,
    7,
    7,
    7,
    8,
    8,
    9,
    7,
    8,
    9,
    8,
    8,
    7,
    8
)

for ($idx = 0; $idx -lt $i0.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0[$idx]
    $ws.Cells.Item($row, 10).Value = $iF[$idx]
}
